$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new labels, shifted from A-D into B-E, with new column A/E
$ws.Range("A1").Value = "HoursTV"
$ws.Range("B1").Value = "Rarely"
$ws.Range("C1").Value = "Sometimes"
$ws.Range("D1").Value = "Often"
$ws.Range("E1").Value = "EveryDay"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108

# Row 2
$ws.Range("A2").Value = "(0,1]"
$ws.Range("B2").Value = 337
$ws.Range("C2").Value = 271
$ws.Range("D2").Value = 83
$ws.Range("E2").Value = 43

# Row 3
$ws.Range("A3").Value = "(1,2]"
$ws.Range("B3").Value = 408
$ws.Range("C3").Value = 331
$ws.Range("D3").Value = 123
$ws.Range("E3").Value = 75

# Row 4
$ws.Range("A4").Value = "(2,3]"
$ws.Range("B4").Value = 283
$ws.Range("C4").Value = 268
$ws.Range("D4").Value = 116
$ws.Range("E4").Value = 79

# Row 5
$ws.Range("A5").Value = "(3,4]"
$ws.Range("B5").Value = 182
$ws.Range("C5").Value = 194
$ws.Range("D5").Value = 64
$ws.Range("E5").Value = 51

# Row 6
$ws.Range("A6").Value = "(4,5]"
$ws.Range("B6").Value = 104
$ws.Range("C6").Value = 109
$ws.Range("D6").Value = 36
$ws.Range("E6").Value = 30

# Row 7
$ws.Range("A7").Value = "(5,10]"
$ws.Range("B7").Value = 124
$ws.Range("C7").Value = 102
$ws.Range("D7").Value = 28
$ws.Range("E7").Value = 46
